$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jrue Holiday"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Boston Celtics"

$ws.Range("A3").Value = "Austin Reaves"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Los Angeles Lakers"

$ws.Range("A4").Value = "Stephen Curry"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Golden State Warriors"

$ws.Range("A5").Value = "Darius Garland"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Cleveland Cavaliers"

$ws.Range("A6").Value = "OG Anunoby"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "New York Knicks"

$ws.Range("A14").Value = "Tyrese Haliburton"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Indiana Pacers"

$ws.Range("A15").Value = "Keegan Murray"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Sacramento Kings"

$ws.Range("A16").Value = "Tyrese Maxey"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Philadelphia 76ers"
